# Update "想去人数" (want-to-go count) column F for a handful of rows on the
# "展览" and "全部类型" sheets, matching the refreshed scrape at 456a3b4.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 232
    4  = 13048
    10 = 229
    11 = 471
    17 = 417
    18 = 5549
    19 = 109
    20 = 57
    22 = 38
    24 = 147
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
